# Updated cryptos list on Mon Aug  7 15:59:06 UTC 2023 with GitHub Actions
#
# This script applies the price/volume refresh (and a small ranking
# re-shuffle among rows 46-50) to the cryptocurrency list on the active
# worksheet, reproducing the target OOXML diff exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while forcing it to stay a text
# value (many "Price" cells look like numbers, e.g. "242.18" or
# "28.952.72") without permanently altering the cell's original style.
function Set-TextValue {
    param($range, $value)
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextValue $ws.Range("D2") '28.952.72'
$ws.Range("E2").Value = '  -0.49%  '
Set-TextValue $ws.Range("D3") '1.825.01'
$ws.Range("E3").Value = '  -0.50%  '
Set-TextValue $ws.Range("D4") '1.004'
$ws.Range("E4").Value = '  +0.11%  '
Set-TextValue $ws.Range("D5") '242.18'
Set-TextValue $ws.Range("D6") '0.6192'
$ws.Range("E6").Value = '  -1.46%  '
Set-TextValue $ws.Range("D7") '1.006'
$ws.Range("E7").Value = '  +0.12%  '
Set-TextValue $ws.Range("D8") '0.07386'
$ws.Range("E8").Value = '  -1.06%  '
Set-TextValue $ws.Range("D9") '0.2904'
$ws.Range("E9").Value = '  -0.56%  '
Set-TextValue $ws.Range("D10") '22.83'
$ws.Range("E10").Value = '  -1.08%  '
Set-TextValue $ws.Range("D11") '0.07673'
$ws.Range("E11").Value = '  -0.77%  '
Set-TextValue $ws.Range("D12") '1.823.44'
$ws.Range("E12").Value = '  -0.51%  '
Set-TextValue $ws.Range("D13") '4.953'
$ws.Range("E13").Value = '  -0.58%  '
Set-TextValue $ws.Range("D14") '0.6660'
$ws.Range("E14").Value = '  -0.34%  '
Set-TextValue $ws.Range("D15") '81.98'
$ws.Range("E15").Value = '  -0.54%  '
Set-TextValue $ws.Range("D16") '0.000009052'
$ws.Range("E16").Value = '  -2.54%  '
Set-TextValue $ws.Range("D17") '5.855'
$ws.Range("E17").Value = '  -2.87%  '
Set-TextValue $ws.Range("D18") '28.961.73'
$ws.Range("E18").Value = '  -0.52%  '
Set-TextValue $ws.Range("D19") '2.079.13'
$ws.Range("E19").Value = '  -0.07%  '
Set-TextValue $ws.Range("D20") '234.45'
$ws.Range("E20").Value = '  +5.19%  '
$ws.Range("E21").Value = '  -0.82%  '
Set-TextValue $ws.Range("D22") '1.005'
$ws.Range("E22").Value = '  +0.04%  '
Set-TextValue $ws.Range("D23") '7.112'
$ws.Range("E23").Value = '  -0.19%  '
$ws.Range("E24").Value = '  +0.14%  '
Set-TextValue $ws.Range("D25") '158.94'
$ws.Range("E25").Value = '  -0.89%  '
Set-TextValue $ws.Range("D26") '0.1405'
$ws.Range("E26").Value = '  +0.64%  '
Set-TextValue $ws.Range("D27") '8.454'
$ws.Range("E27").Value = '  -0.53%  '
Set-TextValue $ws.Range("D28") '17.68'
$ws.Range("E28").Value = '  -1.40%  '
Set-TextValue $ws.Range("D29") '1.488'
$ws.Range("E29").Value = '  -0.97%  '
Set-TextValue $ws.Range("D30") '4.097'
$ws.Range("E30").Value = '  +0.79%  '
Set-TextValue $ws.Range("D31") '0.05529'
$ws.Range("E31").Value = '  -4.43%  '
Set-TextValue $ws.Range("D32") '4.099'
$ws.Range("E32").Value = '  -1.25%  '
$ws.Range("E33").Value = '  +0.32%  '
$ws.Range("E34").Value = '  -1.03%  '
Set-TextValue $ws.Range("D35") '0.7355'
$ws.Range("E35").Value = '  -1.88%  '
Set-TextValue $ws.Range("D36") '1.132'
Set-TextValue $ws.Range("D37") '2.642'
$ws.Range("E37").Value = '  -0.20%  '
Set-TextValue $ws.Range("D38") '2.814'
$ws.Range("E38").Value = '  +2.10%  '
$ws.Range("E39").Value = '  -1.38%  '
Set-TextValue $ws.Range("D40") '1.197.92'
$ws.Range("E40").Value = '  -2.61%  '
Set-TextValue $ws.Range("D41") '6.366'
$ws.Range("E41").Value = '  -2.98%  '
Set-TextValue $ws.Range("D42") '0.9049'
$ws.Range("E42").Value = '  +1.18%  '
$ws.Range("E43").Value = '  +0.07%  '
Set-TextValue $ws.Range("D44") '100.61'
$ws.Range("E44").Value = '  -1.32%  '
Set-TextValue $ws.Range("D45") '1.976.13'
$ws.Range("E45").Value = '  -0.15%  '
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws.Range("D46") '64.36'
$ws.Range("E46").Value = '  -1.95%  '
$ws.Range("B47").Value = 'Mantle'
$ws.Range("C47").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue $ws.Range("D47") '0.5105'
$ws.Range("E47").Value = '  +0.14%  '
$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue $ws.Range("D48") '0.00000000120'
$ws.Range("E48").Value = '  -4.55%  '
$ws.Range("B49").Value = 'TheSandbox'
$ws.Range("C49").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue $ws.Range("D49") '0.4004'
$ws.Range("E49").Value = '  -1.36%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range("D50") '9.031'
$ws.Range("E50").Value = '  +0.47%  '
Set-TextValue $ws.Range("D51") '0.05801'
$ws.Range("E51").Value = '  -0.48%  '